$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 975
$ws.Range("I43").Value = 800
$ws.Range("J43").Value = 1033.3334
$ws.Range("K43").Value = 800
$ws.Range("L43").Value = 1033.3334
$ws.Range("M43").Value = -731
$ws.Range("N43").Value = -1171.3334
$ws.Range("H62").Value = 9270.286
$ws.Range("I62").Value = 7500
$ws.Range("J62").Value = 9978.4
$ws.Range("K62").Value = 7500
$ws.Range("L62").Value = 9978.4
$ws.Range("M62").Value = -6876
$ws.Range("N62").Value = -11226.4
$ws.Range("H65").Value = 9270.286
$ws.Range("I65").Value = 7500
$ws.Range("J65").Value = 9978.4
$ws.Range("K65").Value = 37500
$ws.Range("L65").Value = 49892
$ws.Range("M65").Value = -34380
$ws.Range("N65").Value = -56132
$ws.Range("H107").Value = 740.88464
$ws.Range("I107").Value = 565.8095
$ws.Range("J107").Value = 1476.2
$ws.Range("K107").Value = 565.8095
$ws.Range("L107").Value = 1476.2
$ws.Range("M107").Value = 1354.1905
$ws.Range("N107").Value = -5316.2
$ws.Range("H111").Value = 3140.5881
$ws.Range("I111").Value = 3690
$ws.Range("K111").Value = 11070
$ws.Range("M111").Value = -8003
$ws.Range("H116").Value = 4844.8423
$ws.Range("I116").Value = 2043.75
$ws.Range("K116").Value = 2043.75
$ws.Range("M116").Value = 1398.25
$ws.Range("H124").Value = 33880
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 33880
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 33880
$ws.Range("N124").Value = -43700
$ws.Range("M124").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1815.3846
$ws.Range("I2").Value = 1550
$ws.Range("J2").Value = 5000
$ws.Range("K2").Value = 1550
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = -1437
$ws.Range("N2").Value = -5226
$ws.Range("H32").Value = 5816.92
$ws.Range("I32").Value = 5583.609
$ws.Range("J32").Value = 8500
$ws.Range("K32").Value = 5583.609
$ws.Range("L32").Value = 8500
$ws.Range("M32").Value = -5296.609
$ws.Range("N32").Value = -9074
$ws.Range("H63").Value = 2842153.8
$ws.Range("I63").Value = 1432.5555
$ws.Range("J63").Value = 15625400
$ws.Range("K63").Value = 1432.5555
$ws.Range("L63").Value = 15625400
$ws.Range("M63").Value = -746.5554999999999
$ws.Range("N63").Value = -15626772
$ws.Range("H66").Value = 2842153.8
$ws.Range("I66").Value = 1432.5555
$ws.Range("J66").Value = 15625400
$ws.Range("K66").Value = 7162.7775
$ws.Range("L66").Value = 78127000
$ws.Range("M66").Value = -3730.7775
$ws.Range("N66").Value = -78133864
$ws.Range("H116").Value = 1815.3846
$ws.Range("I116").Value = 1550
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 1550
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = 744
$ws.Range("N116").Value = -9588
$ws.Range("H132").Value = 30779.945
$ws.Range("I132").Value = 2720.5833
$ws.Range("J132").Value = 86898.664
$ws.Range("K132").Value = 8161.749899999999
$ws.Range("L132").Value = 260695.992
$ws.Range("M132").Value = -5631.749899999999
$ws.Range("N132").Value = -265755.992

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1815.3846
$ws.Range("I3").Value = 1550
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 1550
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = -1436
$ws.Range("N3").Value = -5228
$ws.Range("H20").Value = 2766.6667
$ws.Range("I20").Value = 2766.6667
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 2766.6667
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -2519.6667
$ws.Range("N20").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3661.5
$ws.Range("I62").Value = 3599.818
$ws.Range("J62").Value = 3797.2
$ws.Range("K62").Value = 3599.818
$ws.Range("L62").Value = 3797.2
$ws.Range("M62").Value = -2975.818
$ws.Range("N62").Value = -5045.2
$ws.Range("H65").Value = 3661.5
$ws.Range("I65").Value = 3599.818
$ws.Range("J65").Value = 3797.2
$ws.Range("K65").Value = 17999.09
$ws.Range("L65").Value = 18986
$ws.Range("M65").Value = -14879.09
$ws.Range("N65").Value = -25226
$ws.Range("H86").Value = 10427504
$ws.Range("I86").Value = 2618.3635
$ws.Range("J86").Value = 33362252
$ws.Range("K86").Value = 2618.3635
$ws.Range("L86").Value = 33362252
$ws.Range("M86").Value = -1495.3635
$ws.Range("N86").Value = -33364498
$ws.Range("H89").Value = 10427504
$ws.Range("I89").Value = 2618.3635
$ws.Range("J89").Value = 33362252
$ws.Range("K89").Value = 13091.8175
$ws.Range("L89").Value = 166811260
$ws.Range("M89").Value = -7475.817499999999
$ws.Range("N89").Value = -166822492
$ws.Range("H132").Value = 4221.1
$ws.Range("I132").Value = 1599.8572
$ws.Range("K132").Value = 4799.571599999999
$ws.Range("M132").Value = -2269.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 223237.17
$ws.Range("I129").Value = 703.75
$ws.Range("J129").Value = 341921.66
$ws.Range("K129").Value = 2111.25
$ws.Range("L129").Value = 1025764.98
$ws.Range("M129").Value = 2888.75
$ws.Range("N129").Value = -1035764.98
$ws.Range("H131").Value = 715.8
$ws.Range("J131").Value = 724.2708
$ws.Range("L131").Value = 2172.8124
$ws.Range("N131").Value = -12252.8124

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2611589.8
$ws.Range("I70").Value = 3334.8125
$ws.Range("J70").Value = 7828100
$ws.Range("K70").Value = 3334.8125
$ws.Range("L70").Value = 7828100
$ws.Range("M70").Value = -3064.8125
$ws.Range("N70").Value = -7828640
$ws.Range("H73").Value = 2611589.8
$ws.Range("I73").Value = 3334.8125
$ws.Range("J73").Value = 7828100
$ws.Range("K73").Value = 3334.8125
$ws.Range("L73").Value = 7828100
$ws.Range("M73").Value = -2398.8125
$ws.Range("N73").Value = -7829972
$ws.Range("H100").Value = 40666.668
$ws.Range("J100").Value = 40666.668
$ws.Range("L100").Value = 40666.668
$ws.Range("N100").Value = -42830.668
$ws.Range("H113").Value = 2435.8948
$ws.Range("I113").Value = 1939.0834
$ws.Range("K113").Value = 1939.0834
$ws.Range("M113").Value = 230.9166
$ws.Range("H122").Value = 4799.8
$ws.Range("I122").Value = 4749.75
$ws.Range("K122").Value = 14249.25
$ws.Range("M122").Value = -11799.25
$ws.Range("H136").Value = 19556.25
$ws.Range("J136").Value = 19556.25
$ws.Range("L136").Value = 58668.75
$ws.Range("N136").Value = -63768.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 1310500.6
$ws.Range("I122").Value = 2804049.2
$ws.Range("J122").Value = 3645.625
$ws.Range("K122").Value = 8412147.600000001
$ws.Range("L122").Value = 10936.875
$ws.Range("M122").Value = -8409697.600000001
$ws.Range("N122").Value = -15836.875
$ws.Range("H132").Value = 2358.15
$ws.Range("I132").Value = 2019.375
$ws.Range("J132").Value = 3713.25
$ws.Range("K132").Value = 6058.125
$ws.Range("L132").Value = 11139.75
$ws.Range("M132").Value = -3528.125
$ws.Range("N132").Value = -16199.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 40000
$ws.Range("J46").Value = 40000
$ws.Range("L46").Value = 40000
$ws.Range("N46").Value = -40462
$ws.Range("H126").Value = 1206.35
$ws.Range("I126").Value = 1230.3889
$ws.Range("K126").Value = 3691.1667
$ws.Range("M126").Value = -1221.1667
$ws.Range("H134").Value = 40000
$ws.Range("J134").Value = 40000
$ws.Range("L134").Value = 120000
$ws.Range("N134").Value = -125070
